# Generate Report for Handoff
# Adds a new handoff entry (da6fb246-739e-42e4-b222-7ad388577717) as row 3
# to the Overview, zh-cn and de-de sheets, mirroring the existing
# 54780db5-304c-4ec0-95e5-c0ffd5f74adc row already present as row 2.

$wb = $excel.ActiveWorkbook

$fileId = "da6fb246-739e-42e4-b222-7ad388577717"
$zhHash = "3d0f7469a98155ceac1e3a7bcd5d25c99444cd79"

$mdName = "$fileId.md"
$zhXlfName = "$fileId.$zhHash.zh-cn.xlf"
$deXlfName = "$fileId.$zhHash.de-de.xlf"

$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/77df2c771df368776fc5dcf21734e6843cb02678/e2e/$mdName"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fbc4abfea99c085cb49035fa45ff075028c894b8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlfName"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/88357aa6711acee491ede9504217029c19b9c60a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlfName"

$handoffDatetime = "2016-28-20 02:28:28"
$zhHandoffDatetime = "2016-03-20 02:28:25"
$deHandoffDatetime = "2016-03-20 02:28:28"

# ---------------------------------------------------------------
# Sheet "Overview" -- summary row for the new file
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $mdName
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $mdUrl, "", "", $mdName)
$wsOverview.Range("A3").Style = "HyperLink"

$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = $handoffDatetime

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A3").Value = $mdName
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $mdUrl, "", "", $mdName)
$wsZh.Range("A3").Style = "HyperLink"

$wsZh.Range("B3").Value = ".md"
$wsZh.Hyperlinks.Add($wsZh.Range("B3"), $mdUrl, "", "", ".md")
$wsZh.Range("B3").Style = "HyperLink"

$wsZh.Range("C3").Value = "Ready for handoff"

$wsZh.Range("D3").Value = $zhXlfName
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), $zhXlfUrl, "", "", $zhXlfName)
$wsZh.Range("D3").Style = "HyperLink"

$wsZh.Range("E3").Value = $zhHandoffDatetime
$wsZh.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZh.Range("H3").Value = "0001-01-01 00:00:00"
$wsZh.Range("I3").Value = "Include"

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A3").Value = $mdName
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $mdUrl, "", "", $mdName)
$wsDe.Range("A3").Style = "HyperLink"

$wsDe.Range("B3").Value = ".md"
$wsDe.Hyperlinks.Add($wsDe.Range("B3"), $mdUrl, "", "", ".md")
$wsDe.Range("B3").Style = "HyperLink"

$wsDe.Range("C3").Value = "Ready for handoff"

$wsDe.Range("D3").Value = $deXlfName
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), $deXlfUrl, "", "", $deXlfName)
$wsDe.Range("D3").Style = "HyperLink"

$wsDe.Range("E3").Value = $deHandoffDatetime
$wsDe.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDe.Range("I3").Value = "Include"
